$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "Abiola Marcus"
$ws.Range("B3").Value = "Thanks for the brilliant sessions"

[void]$ws.Range("B3").Select()
